$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20: same date as row 19 (2024-08-20), 2 hours, new task description
$ws.Range("A20").NumberFormat = $ws.Range("A19").NumberFormat
$ws.Range("A20").Value = "2024-08-20"

$ws.Range("B20").Value = 2

$ws.Range("C20").Value = "Terminer exercice 2.3 et 2.4"

$ws.Range("B26").Select()
